$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new color values for rows 2-6 (column C), and ux_id values for rows 3-6 (column A)
$ws.Range("C2").Value = "red"

$ws.Range("A3").Value = 1
$ws.Range("C3").Value = "blue"

$ws.Range("A4").Value = 2
$ws.Range("C4").Value = "black"

$ws.Range("A5").Value = 3
$ws.Range("C5").Value = "white"

$ws.Range("A6").Value = 4
$ws.Range("C6").Value = "green"

# Resize the Excel table (ListObject) to cover the new data range
$tbl = $ws.ListObjects.Item("Tabela1")
$tbl.Resize($ws.Range("A1:C6"))

# Update the selected cell on the sheet
$ws.Range("F8").Select()
